$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.086.71"
Set-TextValue "E2" "  +11.46%  "
Set-TextValue "D3" "1.818.33"
Set-TextValue "E3" "  +8.41%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "228.28"
Set-TextValue "E5" "  +3.67%  "
Set-TextValue "D6" "0.546"
Set-TextValue "E6" "  +2.83%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "31.37"
Set-TextValue "E8" "  +4.36%  "
Set-TextValue "D9" "47.06"
Set-TextValue "E9" "  +6.19%  "
Set-TextValue "E10" "  +6.28%  "
Set-TextValue "D11" "0.0667"
Set-TextValue "E11" "  +4.89%  "
Set-TextValue "D12" "0.0929"
Set-TextValue "E12" "  +2.41%  "
Set-TextValue "D13" "2.076.03"
Set-TextValue "E13" "  +8.13%  "
Set-TextValue "D14" "1.826.30"
Set-TextValue "E14" "  +9.10%  "
Set-TextValue "D15" "0.642"
Set-TextValue "E15" "  +4.09%  "
Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "10.36"
Set-TextValue "E16" "  +0.82%  "
Set-TextValue "B17" "WrappedBTC"
Set-TextValue "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "34.084.60"
Set-TextValue "E17" "  +11.41%  "
Set-TextValue "D18" "4.27"
Set-TextValue "E18" "  +7.03%  "
Set-TextValue "D19" "69.47"
Set-TextValue "E19" "  +4.51%  "
Set-TextValue "D20" "258.18"
Set-TextValue "E20" "  +5.15%  "
Set-TextValue "D21" "0.0₃0749"
Set-TextValue "E21" "  +3.42%  "
Set-TextValue "D22" "0.998"
Set-TextValue "E22" "  -0.11%  "
Set-TextValue "D23" "10.51"
Set-TextValue "E23" "  +4.51%  "
Set-TextValue "E24" "  +1.34%  "
Set-TextValue "E25" "  +1.66%  "
Set-TextValue "D26" "157.96"
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "D27" "16.55"
Set-TextValue "E27" "  +3.65%  "
Set-TextValue "D28" "7.15"
Set-TextValue "E28" "  +6.66%  "
Set-TextValue "D29" "0.115"
Set-TextValue "E29" "  +1.76%  "
Set-TextValue "D30" "0.999"
Set-TextValue "E30" "  -0.14%  "
Set-TextValue "E31" "  +10.67%  "
Set-TextValue "D32" "0.0513"
Set-TextValue "E32" "  +3.25%  "
Set-TextValue "E33" "  +4.75%  "
Set-TextValue "E34" "  +6.61%  "
Set-TextValue "B35" "Maker"
Set-TextValue "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D35" "1.542.31"
Set-TextValue "E35" "  +1.94%  "
Set-TextValue "B36" "LidoDAOToken"
Set-TextValue "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D36" "1.81"
Set-TextValue "E36" "  +2.18%  "
Set-TextValue "B37" "TrustWalletToken"
Set-TextValue "C37" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D37" "1.08"
Set-TextValue "E37" "  +5.10%  "
Set-TextValue "B38" "Aave"
Set-TextValue "C38" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D38" "84.95"
Set-TextValue "E38" "  +0.71%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0187"
Set-TextValue "E39" "  +4.65%  "
Set-TextValue "B40" "ImmutableX"
Set-TextValue "C40" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D40" "0.623"
Set-TextValue "E40" "  +3.50%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.82"
Set-TextValue "E41" "  +3.92%  "
Set-TextValue "B42" "HuobiToken"
Set-TextValue "C42" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D42" "2.34"
Set-TextValue "E42" "  +1.48%  "
Set-TextValue "B43" "ARBITRUM"
Set-TextValue "C43" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D43" "0.915"
Set-TextValue "E43" "  +8.73%  "
Set-TextValue "B44" "RenderToken"
Set-TextValue "C44" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "2.15"
Set-TextValue "E44" "  +8.21%  "
Set-TextValue "B45" "Kaspa"
Set-TextValue "C45" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D45" "0.0520"
Set-TextValue "E45" "  +4.33%  "
Set-TextValue "B46" "WEMIXToken"
Set-TextValue "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "1.07"
Set-TextValue "E46" "  +4.52%  "
Set-TextValue "B47" "RocketPoolETH"
Set-TextValue "C47" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D47" "1.976.96"
Set-TextValue "E47" "  +8.91%  "
Set-TextValue "B48" "FraxShare"
Set-TextValue "C48" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D48" "5.72"
Set-TextValue "E48" "  +1.84%  "
Set-TextValue "B49" "PaxDollar"
Set-TextValue "C49" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D49" "0.998"
Set-TextValue "E49" "  -0.13%  "
Set-TextValue "B50" "BitcoinSV"
Set-TextValue "C50" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D50" "52.67"
Set-TextValue "E50" "  +2.14%  "
Set-TextValue "B51" "InjectiveProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "11.65"
Set-TextValue "E51" "  +19.17%  "
